$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 3736
$ws.Range("L3").Value = 3895
$ws.Range("K4").Value = 1772
$ws.Range("L4").Value = 967
$ws.Range("L6").Value = 3385
$ws.Range("K7").Value = 27566
$ws.Range("L7").Value = 12214

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 233
$ws.Range("L3").Value = 264
$ws.Range("L7").Value = 793

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 158
$ws.Range("L3").Value = 185
$ws.Range("L7").Value = 571

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 135
$ws.Range("L3").Value = 143
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 444

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 83
$ws.Range("L7").Value = 231

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 85
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 101
$ws.Range("L6").Value = 98
$ws.Range("L8").Value = 793
$ws.Range("L9").Value = 74
$ws.Range("L10").Value = 79
$ws.Range("L11").Value = 205
$ws.Range("L19").Value = 344
$ws.Range("L20").Value = 306
$ws.Range("L21").Value = 36
$ws.Range("L23").Value = 131
$ws.Range("L25").Value = 65
$ws.Range("L27").Value = 112
$ws.Range("L29").Value = 672
$ws.Range("K31").Value = 328
$ws.Range("L31").Value = 120
$ws.Range("L33").Value = 571
$ws.Range("L34").Value = 75
$ws.Range("L37").Value = 444
$ws.Range("L42").Value = 389
$ws.Range("L43").Value = 91
$ws.Range("L47").Value = 89
$ws.Range("L51").Value = 150
$ws.Range("L54").Value = 253
$ws.Range("L55").Value = 113
$ws.Range("L63").Value = 44
$ws.Range("L65").Value = 231
$ws.Range("L66").Value = 32
$ws.Range("L67").Value = 437
$ws.Range("L72").Value = 54
$ws.Range("L77").Value = 79
$ws.Range("L78").Value = 153
$ws.Range("L79").Value = 318
$ws.Range("L85").Value = 640
$ws.Range("L88").Value = 136
$ws.Range("L89").Value = 173
$ws.Range("L90").Value = 117
$ws.Range("L93").Value = 65
$ws.Range("L96").Value = 123
$ws.Range("L99").Value = 202
$ws.Range("K101").Value = 27566
$ws.Range("L101").Value = 12214

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 43
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 328
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 125
$ws.Range("L3").Value = 166
$ws.Range("L6").Value = 102
$ws.Range("L7").Value = 437

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 61
$ws.Range("L7").Value = 253

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 207
$ws.Range("L3").Value = 253
$ws.Range("L6").Value = 169
$ws.Range("L7").Value = 672

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L2").Value = 22
$ws.Range("L6").Value = 73

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 105
$ws.Range("L7").Value = 344

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 117
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 389

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 35
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 115
$ws.Range("L7").Value = 318

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 98
$ws.Range("L7").Value = 306

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 134
$ws.Range("L6").Value = 112

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 32
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 205

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 173

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L2").Value = 26
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L4").Value = 44
$ws.Range("L6").Value = 136
$ws.Range("L7").Value = 640

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 79
